$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 7, shifting existing rows 7-14 down to 8-15
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the additional NAICS code row
$ws.Cells.Item(7, 1).Value = 333293
$ws.Cells.Item(7, 2).Value = 333244
$ws.Cells.Item(7, 3).Value = 333244

# Update the _FilterDatabase defined name so it covers the expanded range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$15"
    }
}

$wb.Save()
